$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 126
$ws.Range("I2").Value = 308
$ws.Range("J2").Value = 1245
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 359
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 232
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 18
$ws.Range("S2").Value = 117
$ws.Range("T2").Value = 213
$ws.Range("U2").Value = 17
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1900
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 9
